$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 4.795003759225949
$ws.Cells.Item(2, 3).Value = 1.450592490471536
$ws.Cells.Item(2, 4).Value = 0.4824921826647142
$ws.Cells.Item(2, 5).Value = 0.1340698875187947
$ws.Cells.Item(2, 7).Value = 0.002696308362112216
$ws.Cells.Item(2, 10).Value = 0.03518418398326428
$ws.Cells.Item(2, 13).Value = 1.485175881493774
$ws.Cells.Item(2, 14).Value = 4.461129605684732

$ws.Cells.Item(3, 2).Value = 4.627622330456859
$ws.Cells.Item(3, 3).Value = 1.397368962143901
$ws.Cells.Item(3, 4).Value = 0.4812269856431186
$ws.Cells.Item(3, 5).Value = 0.1343454290983743
$ws.Cells.Item(3, 7).Value = 0.002705597747148588
$ws.Cells.Item(3, 10).Value = 0.03532249000028198
$ws.Cells.Item(3, 13).Value = 1.448146701223109
$ws.Cells.Item(3, 14).Value = 4.379294358883243

$ws.Cells.Item(4, 2).Value = 4.52825048008998
$ws.Cells.Item(4, 3).Value = 1.365745384780325
$ws.Cells.Item(4, 4).Value = 0.480689273891528
$ws.Cells.Item(4, 5).Value = 0.1345688355842167
$ws.Cells.Item(4, 7).Value = 0.00271158652076592
$ws.Cells.Item(4, 10).Value = 0.03541221871037514
$ws.Cells.Item(4, 13).Value = 1.426416080122152
$ws.Cells.Item(4, 14).Value = 4.329469525785044

$ws.Cells.Item(5, 2).Value = 4.488600355270364
$ws.Cells.Item(5, 3).Value = 1.353120391228799
$ws.Cells.Item(5, 4).Value = 0.480530062937234
$ws.Cells.Item(5, 5).Value = 0.1346734897746735
$ws.Cells.Item(5, 7).Value = 0.002714098995526167
$ws.Cells.Item(5, 10).Value = 0.03544999709572672
$ws.Cells.Item(5, 13).Value = 1.41781095423562
$ws.Cells.Item(5, 14).Value = 4.309267987235984

$ws.Cells.Item(6, 2).Value = 4.48206719831677
$ws.Cells.Item(6, 3).Value = 1.351039732183096
$ws.Cells.Item(6, 4).Value = 0.4805072383162212
$ws.Cells.Item(6, 5).Value = 0.1346916892156536
$ws.Cells.Item(6, 7).Value = 0.002714520547908529
$ws.Cells.Item(6, 10).Value = 0.03545634356387417
$ws.Cells.Item(6, 13).Value = 1.416397122536836
$ws.Cells.Item(6, 14).Value = 4.305919619300852

$ws.Cells.Item(7, 2).Value = 4.527712338126889
$ws.Cells.Item(7, 3).Value = 1.365574064143175
$ws.Cells.Item(7, 4).Value = 0.4806868844212033
$ws.Cells.Item(7, 5).Value = 0.1345701918892992
$ws.Cells.Item(7, 7).Value = 0.002711620113030445
$ws.Cells.Item(7, 10).Value = 0.03541272328508871
$ws.Cells.Item(7, 13).Value = 1.426299018334149
$ws.Cells.Item(7, 14).Value = 4.329196670689839

$ws.Cells.Item(8, 2).Value = 4.73657830582124
$ws.Cells.Item(8, 3).Value = 1.432019635598351
$ws.Cells.Item(8, 4).Value = 0.4820061768687935
$ws.Cells.Item(8, 5).Value = 0.1341536273455208
$ws.Cells.Item(8, 7).Value = 0.002699452373603709
$ws.Cells.Item(8, 10).Value = 0.03523087643617018
$ws.Cells.Item(8, 13).Value = 1.472197973030859
$ws.Cells.Item(8, 14).Value = 4.432822883087908

$ws.Cells.Item(9, 2).Value = 5.173659256305768
$ws.Cells.Item(9, 3).Value = 1.570872777211321
$ws.Cells.Item(9, 4).Value = 0.486501578570028
$ws.Cells.Item(9, 5).Value = 0.1337680673087327
$ws.Cells.Item(9, 7).Value = 0.00267783849818392
$ws.Cells.Item(9, 10).Value = 0.03491223891863626
$ws.Cells.Item(9, 13).Value = 1.570307319424032
$ws.Cells.Item(9, 14).Value = 4.639563630134944

$ws.Cells.Item(10, 2).Value = 5.51230353787787
$ws.Cells.Item(10, 3).Value = 1.678362692537519
$ws.Cells.Item(10, 4).Value = 0.4909843840962651
$ws.Cells.Item(10, 5).Value = 0.1337494874012144
$ws.Cells.Item(10, 7).Value = 0.002663307861652375
$ws.Cells.Item(10, 10).Value = 0.03470101781795343
$ws.Cells.Item(10, 13).Value = 1.647511893106028
$ws.Cells.Item(10, 14).Value = 4.793872612221605

$ws.Cells.Item(11, 2).Value = 5.670338029825984
$ws.Cells.Item(11, 3).Value = 1.728510176648967
$ws.Cells.Item(11, 4).Value = 0.4932838311952281
$ws.Cells.Item(11, 5).Value = 0.1337989303141001
$ws.Cells.Item(11, 7).Value = 0.00265698594928799
$ws.Cells.Item(11, 10).Value = 0.03460984054731497
$ws.Cells.Item(11, 13).Value = 1.683788719261287
$ws.Cells.Item(11, 14).Value = 4.864655525165972

$ws.Cells.Item(12, 2).Value = 5.73076813779835
$ws.Cells.Item(12, 3).Value = 1.747684160629262
$ws.Cells.Item(12, 4).Value = 0.4941923014936833
$ws.Cells.Item(12, 5).Value = 0.1338260130399505
$ws.Cells.Item(12, 7).Value = 0.002654633088966148
$ws.Cells.Item(12, 10).Value = 0.03457601556881862
$ws.Cells.Item(12, 13).Value = 1.697695409855754
$ws.Cells.Item(12, 14).Value = 4.891548177457878

$ws.Cells.Item(13, 2).Value = 5.717727163815994
$ws.Cells.Item(13, 3).Value = 1.743546433772792
$ws.Cells.Item(13, 4).Value = 0.4939949636062408
$ws.Cells.Item(13, 5).Value = 0.1338198079166375
$ws.Cells.Item(13, 7).Value = 0.00265513799622157
$ws.Cells.Item(13, 10).Value = 0.03458326922966037
$ws.Cells.Item(13, 13).Value = 1.694692767118909
$ws.Cells.Item(13, 14).Value = 4.885752349383324

$ws.Cells.Item(14, 2).Value = 5.675297833207992
$ws.Cells.Item(14, 3).Value = 1.730083913537896
$ws.Cells.Item(14, 4).Value = 0.4933578138528389
$ws.Cells.Item(14, 5).Value = 0.1338009906792301
$ws.Cells.Item(14, 7).Value = 0.002656791555852989
$ws.Cells.Item(14, 10).Value = 0.03460704370234957
$ws.Cells.Item(14, 13).Value = 1.68492941646916
$ws.Cells.Item(14, 14).Value = 4.866866197305171

$ws.Cells.Item(15, 2).Value = 5.649385335413967
$ws.Cells.Item(15, 3).Value = 1.721861857974659
$ws.Cells.Item(15, 4).Value = 0.4929724626598642
$ws.Cells.Item(15, 5).Value = 0.1337905543251594
$ws.Cells.Item(15, 7).Value = 0.002657809752975945
$ws.Cells.Item(15, 10).Value = 0.03462169754129896
$ws.Cells.Item(15, 13).Value = 1.678971245575738
$ws.Cells.Item(15, 14).Value = 4.855309563932735

$ws.Cells.Item(16, 2).Value = 5.502056877434143
$ws.Cells.Item(16, 3).Value = 1.675110967261617
$ws.Cells.Item(16, 4).Value = 0.4908393713149337
$ws.Cells.Item(16, 5).Value = 0.1337474240938015
$ws.Cells.Item(16, 7).Value = 0.002663726784313048
$ws.Cells.Item(16, 10).Value = 0.03470707489772096
$ws.Cells.Item(16, 13).Value = 1.645164651854401
$ws.Cells.Item(16, 14).Value = 4.78925895892155

$ws.Cells.Item(17, 2).Value = 5.412705077388523
$ws.Cells.Item(17, 3).Value = 1.646754071902421
$ws.Cells.Item(17, 4).Value = 0.48959764359536
$ws.Cells.Item(17, 5).Value = 0.1337358173504626
$ws.Cells.Item(17, 7).Value = 0.002667430262297841
$ws.Cells.Item(17, 10).Value = 0.03476070546715881
$ws.Cells.Item(17, 13).Value = 1.624723747690794
$ws.Cells.Item(17, 14).Value = 4.748892419392917

$ws.Cells.Item(18, 2).Value = 5.36168642819905
$ws.Cells.Item(18, 3).Value = 1.630561293941639
$ws.Cells.Item(18, 4).Value = 0.4889079066673361
$ws.Cells.Item(18, 5).Value = 0.1337345886275472
$ws.Cells.Item(18, 7).Value = 0.002669587546062109
$ws.Cells.Item(18, 10).Value = 0.03479201457549408
$ws.Cells.Item(18, 13).Value = 1.613075298456309
$ws.Cells.Item(18, 14).Value = 4.725729641411334

$ws.Cells.Item(19, 2).Value = 5.344476258191889
$ws.Cells.Item(19, 3).Value = 1.625098718895686
$ws.Cells.Item(19, 4).Value = 0.4886785676341248
$ws.Cells.Item(19, 5).Value = 0.1337351071388184
$ws.Cells.Item(19, 7).Value = 0.002670322636827328
$ws.Cells.Item(19, 10).Value = 0.03480269480871812
$ws.Cells.Item(19, 13).Value = 1.609149893029596
$ws.Cells.Item(19, 14).Value = 4.717896437567703

$ws.Cells.Item(20, 2).Value = 5.422177913364976
$ws.Cells.Item(20, 3).Value = 1.649760533491644
$ws.Cells.Item(20, 4).Value = 0.4897272927101426
$ws.Cells.Item(20, 5).Value = 0.1337364889029935
$ws.Cells.Item(20, 7).Value = 0.002667033214151357
$ws.Cells.Item(20, 10).Value = 0.034754948591603
$ws.Cells.Item(20, 13).Value = 1.62688845239046
$ws.Cells.Item(20, 14).Value = 4.753183784585644

$ws.Cells.Item(21, 2).Value = 5.687744349651211
$ws.Cells.Item(21, 3).Value = 1.734033145845842
$ws.Cells.Item(21, 4).Value = 0.4935439340005274
$ws.Cells.Item(21, 5).Value = 0.1338062905911706
$ws.Cells.Item(21, 7).Value = 0.002656304752086981
$ws.Cells.Item(21, 10).Value = 0.0346000415456178
$ws.Cells.Item(21, 13).Value = 1.687792524559967
$ws.Cells.Item(21, 14).Value = 4.872411076061155

$ws.Cells.Item(22, 2).Value = 5.864728637707117
$ws.Cells.Item(22, 3).Value = 1.790186130527673
$ws.Cells.Item(22, 4).Value = 0.4962583381767303
$ws.Cells.Item(22, 5).Value = 0.1339006545985804
$ws.Cells.Item(22, 7).Value = 0.002649532557401037
$ws.Cells.Item(22, 10).Value = 0.03450289008358887
$ws.Cells.Item(22, 13).Value = 1.728586061915337
$ws.Cells.Item(22, 14).Value = 4.950852065045012

$ws.Cells.Item(23, 2).Value = 5.769951331109155
$ws.Cells.Item(23, 3).Value = 1.76011625281933
$ws.Cells.Item(23, 4).Value = 0.4947893752182893
$ws.Cells.Item(23, 5).Value = 0.1338458186935121
$ws.Cells.Item(23, 7).Value = 0.002653125200762111
$ws.Cells.Item(23, 10).Value = 0.03455436877564821
$ws.Cells.Item(23, 13).Value = 1.706722187228209
$ws.Cells.Item(23, 14).Value = 4.908937698001239

$ws.Cells.Item(24, 2).Value = 5.417894152255144
$ws.Cells.Item(24, 3).Value = 1.648400969985687
$ws.Cells.Item(24, 4).Value = 0.4896686031496387
$ws.Cells.Item(24, 5).Value = 0.1337361683394107
$ws.Cells.Item(24, 7).Value = 0.002667212632007323
$ws.Cells.Item(24, 10).Value = 0.03475754979122669
$ws.Cells.Item(24, 13).Value = 1.625909467883773
$ws.Cells.Item(24, 14).Value = 4.751243520150013

$ws.Cells.Item(25, 2).Value = 5.052396500170289
$ws.Cells.Item(25, 3).Value = 1.532367147651883
$ws.Cells.Item(25, 4).Value = 0.4850794731755173
$ws.Cells.Item(25, 5).Value = 0.133826021480548
$ws.Cells.Item(25, 7).Value = 0.002683447207246416
$ws.Cells.Item(25, 10).Value = 0.03499440164945966
$ws.Cells.Item(25, 13).Value = 1.542880034890302
$ws.Cells.Item(25, 14).Value = 4.583229783554913
